$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Simple per-row D/E updates ---
Set-TextValue 'D2' '29.284.45'
Set-TextValue 'E2' '  -0.77%  '

Set-TextValue 'D3' '1.869.81'
Set-TextValue 'E3' '  -0.46%  '

Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  -0.14%  '

Set-TextValue 'D5' '0.7111'
Set-TextValue 'E5' '  -1.13%  '

Set-TextValue 'D6' '241.82'
Set-TextValue 'E6' '  +0.00%  '

Set-TextValue 'E7' '  -0.15%  '

Set-TextValue 'D8' '0.3111'
Set-TextValue 'E8' '  +0.35%  '

Set-TextValue 'D9' '0.07702'
Set-TextValue 'E9' '  -2.56%  '

Set-TextValue 'E10' '  -3.09%  '

Set-TextValue 'D11' '0.08398'

Set-TextValue 'D12' '1.883.75'
Set-TextValue 'E12' '  -0.34%  '

Set-TextValue 'D13' '5.223'
Set-TextValue 'E13' '  -1.13%  '

Set-TextValue 'D14' '0.7121'
Set-TextValue 'E14' '  -2.38%  '

Set-TextValue 'E15' '  -0.26%  '

Set-TextValue 'D16' '29.282.71'
Set-TextValue 'E16' '  -0.82%  '

Set-TextValue 'D17' '0.000008150'
Set-TextValue 'E17' '  +3.68%  '

Set-TextValue 'D18' '5.937'
Set-TextValue 'E18' '  +0.50%  '

Set-TextValue 'D19' '243.59'
Set-TextValue 'E19' '  -0.95%  '

Set-TextValue 'D20' '2.125.13'
Set-TextValue 'E20' '  -0.81%  '

Set-TextValue 'D21' '13.13'
Set-TextValue 'E21' '  -1.42%  '

Set-TextValue 'D22' '0.9998'
Set-TextValue 'E22' '  -0.24%  '

Set-TextValue 'D23' '7.879'
Set-TextValue 'E23' '  -2.54%  '

Set-TextValue 'E24' '  -0.16%  '

Set-TextValue 'D25' '0.1623'
Set-TextValue 'E25' '  -0.04%  '

Set-TextValue 'D26' '164.28'
Set-TextValue 'E26' '  +0.47%  '

Set-TextValue 'E27' '  -0.47%  '

Set-TextValue 'E28' '  +0.94%  '

Set-TextValue 'D29' '1.510'
Set-TextValue 'E29' '  +1.00%  '

Set-TextValue 'D30' '4.403'
Set-TextValue 'E30' '  +0.11%  '

Set-TextValue 'D31' '1.308'
Set-TextValue 'E31' '  -3.42%  '

Set-TextValue 'D32' '4.282'
Set-TextValue 'E32' '  +4.14%  '

Set-TextValue 'D33' '0.05173'
Set-TextValue 'E33' '  -0.75%  '

Set-TextValue 'E36' '  -2.44%  '

Set-TextValue 'D37' '2.685'
Set-TextValue 'E37' '  +0.27%  '

Set-TextValue 'E38' '  -0.65%  '

Set-TextValue 'D39' '2.712'
Set-TextValue 'E39' '  +0.26%  '

Set-TextValue 'D40' '1.157.77'
Set-TextValue 'E40' '  -3.95%  '

Set-TextValue 'D41' '6.398'
Set-TextValue 'E41' '  +3.51%  '

Set-TextValue 'E44' '  -0.16%  '

Set-TextValue 'D45' '103.28'
Set-TextValue 'E45' '  +0.81%  '

Set-TextValue 'D46' '2.020.39'
Set-TextValue 'E46' '  -0.21%  '

Set-TextValue 'D49' '9.391'
Set-TextValue 'E49' '  +0.76%  '

Set-TextValue 'D50' '0.4296'
Set-TextValue 'E50' '  -0.75%  '

Set-TextValue 'D51' '7.046'
Set-TextValue 'E51' '  -0.59%  '

# --- Row swaps (coin rows reordered) ---
# Swap rows 34 and 35
Set-TextValue 'B34' 'LidoDAOToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue 'D34' '1.915'
Set-TextValue 'E34' '  -1.75%  '
Set-TextValue 'B35' 'ImmutableX'
Set-TextValue 'C35' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D35' '0.7734'
Set-TextValue 'E35' '  +6.25%  '

# Swap rows 42 and 43
Set-TextValue 'B42' 'TrustWalletToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '0.8915'
Set-TextValue 'E42' '  -2.01%  '
Set-TextValue 'B43' 'Aave'
Set-TextValue 'C43' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D43' '73.24'
Set-TextValue 'E43' '  -0.39%  '

# Swap rows 47 and 48
Set-TextValue 'B47' 'RenderToken'
Set-TextValue 'C47' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D47' '1.796'
Set-TextValue 'E47' '  -0.26%  '
Set-TextValue 'B48' 'Mantle'
Set-TextValue 'C48' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D48' '0.5191'
Set-TextValue 'E48' '  -1.93%  '

